$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 646:647, shifting all rows from 646 onward down by 2
$ws.Range("646:647").Insert()

# New row 646 data
$ws.Range("A646").Value = 10
$ws.Range("B646").Value = "Vega Modelo de Temuco"
$ws.Range("C646").Value = "La Araucanía"
$ws.Range("D646").Value = 44769
$ws.Range("E646").Value = 9
$ws.Range("F646").Value = 100112003
$ws.Range("G646").Value = "Ajo"
$ws.Range("H646").Value = "Chino"
$ws.Range("I646").Value = "Primera"
$ws.Range("J646").Value = 300
$ws.Range("K646").Value = 30000
$ws.Range("L646").Value = 33000
$ws.Range("M646").Value = 31000
$ws.Range("N646").Value = "$/caja 10 kilos"
$ws.Range("O646").Value = "China"
$ws.Range("P646").Value = 3100
$ws.Range("Q646").Value = 10
$ws.Range("R646").Value = "Hortaliza"

# New row 647 data
$ws.Range("A647").Value = 10
$ws.Range("B647").Value = "Vega Modelo de Temuco"
$ws.Range("C647").Value = "La Araucanía"
$ws.Range("D647").Value = 44769
$ws.Range("E647").Value = 9
$ws.Range("F647").Value = 100112003
$ws.Range("G647").Value = "Ajo"
$ws.Range("H647").Value = "Chino"
$ws.Range("I647").Value = "Primera"
$ws.Range("J647").Value = 100
$ws.Range("K647").Value = 35000
$ws.Range("L647").Value = 35000
$ws.Range("M647").Value = 35000
$ws.Range("N647").Value = "$/malla 10 kilos"
$ws.Range("O647").Value = "China"
$ws.Range("P647").Value = 3500
$ws.Range("Q647").Value = 10
$ws.Range("R647").Value = "Hortaliza"
